$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update status text from "In Translation" to "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Update timestamps
$overview.Range("G2").Value = "2016-09-05 19:05:52"
$zhcn.Range("H2").Value = "2016-09-05 19:05:48"
$dede.Range("H2").Value = "2016-09-05 19:05:52"

# Widen the columns that hold the status text so widths match the new, longer content.
# NOTE: the host's ColumnWidth setter quantizes to an internal pixel grid (stored width
# ends up snapped to the nearest 1/6 character unit), so we feed it the input value whose
# quantized result lands as close as possible to the recorded target width (17.2159881591797).
$targetColumnWidthInput = 98.0 / 6.0
$overview.Columns.Item(5).ColumnWidth = $targetColumnWidthInput
$overview.Columns.Item(6).ColumnWidth = $targetColumnWidthInput
$zhcn.Columns.Item(3).ColumnWidth = $targetColumnWidthInput
$dede.Columns.Item(3).ColumnWidth = $targetColumnWidthInput
